# Table S_Predictors.xlsx - add a "Notes" column (F) to the predictors sheet
# with supplementary notes for a handful of predictor rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New header for column F
$ws.Range("F1").Value = 'Notes'
$ws.Range("F1").Font.Bold = $true

# Notes for specific predictor rows
$ws.Range("F16").Value = 'TPI (Topographic Position Index) is the difference between the value of a central cell and the mean value of its surrounding cells within 1 km window. - raster package'
$ws.Range("F17").Value = 'TPI (Topographic Position Index) is the difference between the value of the central cell and the mean value of its surrounding cells within 250m window. - raster package'
$ws.Range("F18").Value = 'TPI (Topographic Position Index) is the difference between the value of the central cell and the mean value of its surrounding cells within 500 m window - raster package'
$ws.Range("F20").Value = 'topographic wetness index represents a theoretical estimation of the accumulation of flow at any point (ie depends on upstream catchment area). Formula: ln(a/tan(beta)); (a = upslope contributing area per unit contour; tan(beta) = local slope angle)'
$ws.Range("F27").Value = 'NBR: Normalised Burn Ratio - identify burn areas and severity of burn'
$ws.Range("F39").Value = 'NDMI: Normalised Difference Moisture Index, for vegetation water content'
$ws.Range("F54").Value = 'From Oregon shapefile of logged areas over last ~ 100 years'

# Move selection/viewport so it matches the refreshed view of the sheet
$ws.Range("F1").Select()
